$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 2957
$ws.Range("F7").Value = 241
$ws.Range("F10").Value = 6955
$ws.Range("F11").Value = 43
$ws.Range("F12").Value = 92
$ws.Range("F13").Value = 360
$ws.Range("F14").Value = 607
$ws.Range("F15").Value = 1500
$ws.Range("F17").Value = 2250
$ws.Range("F18").Value = 1502
$ws.Range("F20").Value = 126
$ws.Range("F22").Value = 136
$ws.Range("F23").Value = 186
$ws.Range("F26").Value = 1746
$ws.Range("F30").Value = 1669
$ws.Range("F33").Value = 590
$ws.Range("F34").Value = 425
$ws.Range("F35").Value = 433
$ws.Range("F36").Value = 21
$ws.Range("F37").Value = 2493
$ws.Range("F38").Value = 2737
$ws.Range("F39").Value = 76
$ws.Range("F40").Value = 16
$ws.Range("F42").Value = 19
$ws.Range("F43").Value = 30
$ws.Range("F44").Value = 322
$ws.Range("F47").Value = 160
$ws.Range("F48").Value = 15

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 170
$ws.Range("F20").Value = 57
$ws.Range("F22").Value = 341
$ws.Range("F23").Value = 482

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 543
$ws.Range("F6").Value = 1697
$ws.Range("F8").Value = 2745
$ws.Range("F9").Value = 1026
$ws.Range("F10").Value = 941
$ws.Range("F12").Value = 278
$ws.Range("F13").Value = 1496
$ws.Range("F14").Value = 7383

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2957
$ws.Range("F5").Value = 241
$ws.Range("F6").Value = 1697
$ws.Range("F7").Value = 2745
$ws.Range("F8").Value = 6955
$ws.Range("F9").Value = 1026
$ws.Range("F10").Value = 43
$ws.Range("F11").Value = 360
$ws.Range("F12").Value = 170
$ws.Range("F13").Value = 1496
$ws.Range("F14").Value = 607
$ws.Range("F15").Value = 1500
$ws.Range("F17").Value = 2250
$ws.Range("F18").Value = 1502
$ws.Range("F19").Value = 126
$ws.Range("F22").Value = 136
$ws.Range("F25").Value = 1746
$ws.Range("F29").Value = 1669
$ws.Range("F32").Value = 590
$ws.Range("F33").Value = 434
$ws.Range("F34").Value = 57
$ws.Range("F36").Value = 482
$ws.Range("F37").Value = 433
$ws.Range("F38").Value = 21
$ws.Range("F39").Value = 2493
$ws.Range("F40").Value = 2737
$ws.Range("F41").Value = 76
$ws.Range("F43").Value = 30
$ws.Range("F44").Value = 322
